{"js": "// Update the worksheet date and all the division problems/answers.\n// Each (oldText -> newText) pair below corresponds 1:1 to a single\n// <w:t> run in the document; every \"old\" string is unique in the\n// document, so a plain case-sensitive search-and-replace is safe and\n// will not accidentally touch a run that was already rewritten by an\n// earlier pair in this same list.\nconst replacements = [\n  [\"2023-07-23 Sunday\", \"2023-07-24 Monday\"],\n  [\"85\u00f79=9, 4\", \"14\u00f77=2, 0\"],\n  [\"70\u00f78=8, 6\", \"21\u00f77=3, 0\"],\n  [\"49\u00f72=24, 1\", \"61\u00f79=6, 7\"],\n  [\"58\u00f77=8, 2\", \"66\u00f74=16, 2\"],\n  [\"28\u00f78=3, 4\", \"21\u00f76=3, 3\"],\n  [\"82\u00f72=41, 0\", \"10\u00f72=5, 0\"],\n  [\"84\u00f78=10, 4\", \"83\u00f78=10, 3\"],\n  [\"46\u00f76=7, 4\", \"63\u00f75=12, 3\"],\n  [\"68\u00f74=17, 0\", \"92\u00f78=11, 4\"],\n  [\"65\u00f79=7, 2\", \"43\u00f77=6, 1\"],\n  [\"86\u00f74=21, 2\", \"30\u00f75=6, 0\"],\n  [\"59\u00f77=8, 3\", \"73\u00f77=10, 3\"],\n  [\"93\u00f74=23, 1\", \"76\u00f75=15, 1\"],\n  [\"22\u00f72=11, 0\", \"81\u00f74=20, 1\"],\n  [\"58\u00f72=29, 0\", \"29\u00f74=7, 1\"],\n  [\"24\u00f75=4, 4\", \"30\u00f73=10, 0\"],\n  [\"85\u00f73=28, 1\", \"94\u00f77=13, 3\"],\n  [\"22\u00f78=2, 6\", \"60\u00f73=20, 0\"],\n  [\"94\u00f73=31, 1\", \"80\u00f75=16, 0\"],\n  [\"35\u00f76=5, 5\", \"33\u00f75=6, 3\"],\n  [\"10\u00f76=1, 4\", \"10\u00f72=5, 0\"],\n  [\"76\u00f74=19, 0\", \"19\u00f79=2, 1\"],\n  [\"42\u00f72=21, 0\", \"24\u00f79=2, 6\"],\n  [\"79\u00f79=8, 7\", \"94\u00f74=23, 2\"],\n  [\"45\u00f73=15, 0\", \"57\u00f75=11, 2\"],\n];\n\nconst body = context.document.body;\n\nfor (const [oldText, newText] of replacements) {\n  const results = body.search(oldText, { matchCase: true, matchWholeWord: false });\n  results.load(\"items\");\n  await context.sync();\n\n  for (const range of results.items) {\n    range.insertText(newText, Word.InsertLocation.replace);\n  }\n  await context.sync();\n}\n", "ps1": "# Update the worksheet date and all the division problems/answers.\n# Each (OldText -> NewText) pair corresponds 1:1 to a single <w:t> run\n# in the document; every \"OldText\" string is unique in the document,\n# so Find/Replace (one match per call, wdReplaceAll) is safe and will\n# not accidentally touch a run that was already rewritten by an\n# earlier pair in this same list.\n\n$d = $word.ActiveDocument\n\n$replacements = @(\n    @(\"2023-07-23 Sunday\", \"2023-07-24 Monday\"),\n    @(\"85\u00f79=9, 4\", \"14\u00f77=2, 0\"),\n    @(\"70\u00f78=8, 6\", \"21\u00f77=3, 0\"),\n    @(\"49\u00f72=24, 1\", \"61\u00f79=6, 7\"),\n    @(\"58\u00f77=8, 2\", \"66\u00f74=16, 2\"),\n    @(\"28\u00f78=3, 4\", \"21\u00f76=3, 3\"),\n    @(\"82\u00f72=41, 0\", \"10\u00f72=5, 0\"),\n    @(\"84\u00f78=10, 4\", \"83\u00f78=10, 3\"),\n    @(\"46\u00f76=7, 4\", \"63\u00f75=12, 3\"),\n    @(\"68\u00f74=17, 0\", \"92\u00f78=11, 4\"),\n    @(\"65\u00f79=7, 2\", \"43\u00f77=6, 1\"),\n    @(\"86\u00f74=21, 2\", \"30\u00f75=6, 0\"),\n    @(\"59\u00f77=8, 3\", \"73\u00f77=10, 3\"),\n    @(\"93\u00f74=23, 1\", \"76\u00f75=15, 1\"),\n    @(\"22\u00f72=11, 0\", \"81\u00f74=20, 1\"),\n    @(\"58\u00f72=29, 0\", \"29\u00f74=7, 1\"),\n    @(\"24\u00f75=4, 4\", \"30\u00f73=10, 0\"),\n    @(\"85\u00f73=28, 1\", \"94\u00f77=13, 3\"),\n    @(\"22\u00f78=2, 6\", \"60\u00f73=20, 0\"),\n    @(\"94\u00f73=31, 1\", \"80\u00f75=16, 0\"),\n    @(\"35\u00f76=5, 5\", \"33\u00f75=6, 3\"),\n    @(\"10\u00f76=1, 4\", \"10\u00f72=5, 0\"),\n    @(\"76\u00f74=19, 0\", \"19\u00f79=2, 1\"),\n    @(\"42\u00f72=21, 0\", \"24\u00f79=2, 6\"),\n    @(\"79\u00f79=8, 7\", \"94\u00f74=23, 2\"),\n    @(\"45\u00f73=15, 0\", \"57\u00f75=11, 2\")\n)\n\nforeach ($pair in $replacements) {\n    $oldText = $pair[0]\n    $newText = $pair[1]\n\n    $find = $d.Content.Find\n    $find.ClearFormatting()\n    $find.Replacement.ClearFormatting()\n    $find.Text = $oldText\n    $find.Replacement.Text = $newText\n    $find.Execute($oldText, $false, $false, $false, $false, $false, $true, 1, $false, $newText, 2)\n}\n"}
